$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 25: new daily report data (28 de Marzo) ---
$ws.Range("A25").Value = 43916
$ws.Range("A25").NumberFormat = "DD/MM/YY"
$ws.Range("B25").Value = 24
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 21
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = 14
$ws.Range("H25").Value = 49
$ws.Range("I25").Value = 938
$ws.Range("J25").Value = 16
$ws.Range("K25").Value = 32
$ws.Range("L25").Value = 144
$ws.Range("M25").Value = 135
$ws.Range("N25").Value = 143
$ws.Range("O25").Value = 22
$ws.Range("P25").Value = 63
$ws.Range("Q25").Value = 2
$ws.Range("R25").Value = 22
$ws.Range("S25").Value = 1610

# --- Row 34: new highlighted / wrapped note band (K34:O34) ---
$hdr = $ws.Range("K34")
$hdr.Font.Bold = $true
$hdr.Interior.Color = 14788918
$hdr.Interior.PatternColor = 16763904
$hdr.WrapText = $true

# propagate the same look to the rest of the band without re-walking
# the whole style chain cell by cell
$hdr.Copy()
$ws.Range("L34:O34").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Selection / scroll state left by the edit ---
$ws.Range("J33:Q56").Select()
